# This document was re-saved by a different OOXML writer, which merges
# adjacent runs that carry identical run formatting into a single <w:r>
# (and, as a side effect, drops the interior <w:proofErr/> spell/grammar
# markers that used to split them). No visible word, sentence or
# formatting actually changed - only how the runs are split. We
# reproduce that by re-typing the affected paragraphs' text in place via
# Find & Replace scoped to each paragraph's own Range, which collapses
# the multiple runs Word previously tracked into one run sharing the
# paragraph/first-run formatting, exactly like the target.

$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    param($ParagraphIndex, $Text)
    $p = $d.Paragraphs($ParagraphIndex)
    $rng = $p.Range
    $rng.Find.ClearFormatting()
    [void]$rng.Find.Execute($Text, $false, $false, $false, $false, $false, $true, 1, $false, $Text, 2)
}

# "Questions for Tree app: (We can try to put images of what these look
# like for each question)" - was split across 3 runs around a
# proofErr gramStart/gramEnd pair ("these look").
Merge-ParagraphRuns 1 "Questions for Tree app: (We can try to put images of what these look like for each question)"

# "Q5A: Do the Leaves have deep SHARP V-shaped Sinuses (indentations
# between each leaf lobe) or SMOOTH Sinuses?" - 5 italic runs merged.
Merge-ParagraphRuns 10 "Q5A: Do the Leaves have deep SHARP V-shaped Sinuses (indentations between each leaf lobe) or SMOOTH Sinuses?"

# "Q2B: Is the leaf SIMPLE or COMPOUND?" - 3 italic runs merged.
Merge-ParagraphRuns 22 "Q2B: Is the leaf SIMPLE or COMPOUND?"

# "Q5B: Are the leaves SERRATED with DEEP LOBES, or are they big, WAXY
# AND GLOSSY?" - 2 italic runs merged (also drops the transient
# lastRenderedPageBreak marker that Word had cached on the first run).
Merge-ParagraphRuns 34 "Q5B: Are the leaves SERRATED with DEEP LOBES, or are they big, WAXY AND GLOSSY?"

# "Q6D: Are the compound leaflets growing INDIVIDUALLY or are there
# SEVERAL" - 4 italic runs ("Q" / "6" / "D:" / " Are ... SEVERAL")
# merged (also drops a cached lastRenderedPageBreak marker).
Merge-ParagraphRuns 74 "Q6D: Are the compound leaflets growing INDIVIDUALLY or are there SEVERAL"

# "Jan Yaro" - bold run split by a spellStart/spellEnd proofErr pair.
# Scope the Find to just the bold name so the preceding plain
# "IDENTIFICATION: " run is left untouched.
$p80 = $d.Paragraphs(80)
$nameRng = $p80.Range
$nameRng.Find.ClearFormatting()
[void]$nameRng.Find.Execute("Jan Yaro", $false, $false, $false, $false, $false, $true, 1, $false, "Jan Yaro", 2)
